# Refresh cryptocurrency price/volume snapshot (scheduled GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.010.02'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '2.912.31'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'371.60"
$ws.Range("E5").Value = '  +4.88%  '
$ws.Range("D6").Value = "'101.94"
$ws.Range("E6").Value = '  -4.84%  '
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = '  -3.37%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -5.49%  '
$ws.Range("D10").Value = "'36.87"
$ws.Range("E10").Value = '  -3.85%  '
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = "'0.0834"
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("D13").Value = "'18.27"
$ws.Range("E13").Value = '  -4.74%  '
$ws.Range("D14").Value = '3.368.35'
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").Value = "'7.35"
$ws.Range("E15").Value = '  -3.55%  '
$ws.Range("D16").Value = '2.911.78'
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("E17").Value = '  -8.62%  '
$ws.Range("D18").Value = '50.936.84'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").Value = "'3.22"
$ws.Range("E19").Value = '  -6.94%  '
$ws.Range("E20").Value = '  -4.11%  '
$ws.Range("D21").Value = "'12.88"
$ws.Range("E21").Value = '  -5.05%  '
$ws.Range("E22").Value = '  -3.32%  '
$ws.Range("D23").Value = "'67.96"
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").Value = "'258.46"
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = '  -5.75%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'25.55"
$ws.Range("E28").Value = '  -4.96%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = "'7.06"
$ws.Range("E29").Value = '  -5.37%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.102"
$ws.Range("E30").Value = '  -5.11%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = "'6.27"
$ws.Range("E31").Value = '  +3.11%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = "'9.84"
$ws.Range("E32").Value = '  -4.49%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = "'2.11"
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = "'51.27"
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").Value = "'34.09"
$ws.Range("E35").Value = '  -5.67%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0421"
$ws.Range("E37").Value = '  -4.08%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = "'2.98"
$ws.Range("E38").Value = '  -7.54%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = "'16.99"
$ws.Range("E39").Value = '  -4.98%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = "'2.59"
$ws.Range("E40").Value = '  -4.64%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = "'1.84"
$ws.Range("E41").Value = '  -6.51%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = "'0.112"
$ws.Range("E42").Value = '  -3.38%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = "'119.32"
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").Value = "'21.83"
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'2.09"
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.020.57'
$ws.Range("E46").Value = '  -4.57%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = "'2.31"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'3.13"
$ws.Range("E48").Value = '  -6.91%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '3.195.05'
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = "'0.235"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = "'0.0308"
$ws.Range("E51").Value = '  -9.40%  '

# Clear the quote-prefix style marker picked up above so formatting matches plain text cells
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
